$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at position 4 (pushes the old "disability persons"
#    row to 5 and the merged Source row to 6), keeping formulas/merges
#    shifted along with it.
# ---------------------------------------------------------------------
$ws.Rows("4:4").Insert()

# ---------------------------------------------------------------------
# 2. New row 4: "family with disabilities Persons " + counts
#    Borrow formatting from row 5 (the old label row, which still has
#    the matching font/fill/top+bottom border) then trim the bottom
#    border off so only the top edge remains.
# ---------------------------------------------------------------------
$ws.Range("A5").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B5:I5").Copy()
$ws.Range("B4:I4").PasteSpecial(-4122)

$ws.Range("A4").Value = "family with disabilities Persons "
$ws.Range("B4").Value = 899
$ws.Range("C4").Value = 847
$ws.Range("D4").Value = 782
$ws.Range("E4").Value = 805
$ws.Range("F4").Value = 804
$ws.Range("G4").Value = 805
$ws.Range("H4").Value = 795
$ws.Range("I4").Value = 782

$ws.Range("A4").Borders.LineStyle = -4142
$ws.Range("A4").Borders.Item(8).LineStyle = 1
$ws.Range("I4").HorizontalAlignment = -4131

$ws.Rows("4:4").RowHeight = 24.75

# ---------------------------------------------------------------------
# 3. Row 5 (the shifted-down former row 4): relabel + new counts, and
#    change its border so only the bottom edge remains.
# ---------------------------------------------------------------------
$ws.Range("A5").Value = "disabilities Persons "
$ws.Range("B5").Value = 978
$ws.Range("C5").Value = 927
$ws.Range("D5").Value = 856
$ws.Range("E5").Value = 877
$ws.Range("F5").Value = 872
$ws.Range("G5").Value = 873
$ws.Range("H5").Value = 860
$ws.Range("I5").Value = 845

$ws.Range("A5").Borders.LineStyle = -4142
$ws.Range("A5").Borders.Item(9).LineStyle = 1
$ws.Range("I5").HorizontalAlignment = -4131
$ws.Range("I5").Borders.Item(9).LineStyle = 1

$ws.Rows("5:5").RowHeight = 21

# ---------------------------------------------------------------------
# 4. Row 6 (the shifted-down Source row): text is unchanged, just keeps
#    moving along with the insert; tidy its height to match the target.
# ---------------------------------------------------------------------
$ws.Rows("6:6").RowHeight = 27.75

# ---------------------------------------------------------------------
# 5. Row 2 loses its explicit custom height (goes back to default).
# ---------------------------------------------------------------------
$ws.Rows("2:2").AutoFit()

# ---------------------------------------------------------------------
# 6. Row 1: new merged title spanning A1:I1, bold + centered + wrapped,
#    with the new wording, and a taller row height.
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "The number of persons with disabilities registered in the Unified database of targeted social assistance program in Kaspi Municipality"
$ws.Range("A1:I1").Merge()
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").WrapText = $true
$ws.Rows("1:1").RowHeight = 51
$ws.Range("A1:I1").Select()

# ---------------------------------------------------------------------
# 7. Column A width tweak to match the new layout.
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 20
